# Generate Report for Handback
#
# The nightly localization-status report is refreshed once a handback
# completes: the "Status" column flips from "Ready for handoff" to
# "Handed back: in sync with en-US" for every locale, and each per-locale
# sheet gets its "Latest Target File" / "Latest Handback File" /
# "Latest Handback DateTime" columns populated (incl. a fresh hyperlink to
# the source .md on the Latest Target File cell). The "…DateTime" /
# "…File" columns are now wider so the new text fits.

$wb = $excel.ActiveWorkbook

$mdDisplay = "7eacbddf-b45d-4045-a9d6-3e4957d931f6.md"
$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/0d7f953a4986a9c33572b2ce1e80a8f6a170725e/e2e/7eacbddf-b45d-4045-a9d6-3e4957d931f6.md"
$handedBackStatus = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR encoding of RGB 6495ED -- matches the workbook's existing HyperLink cell style

# ---------------------------------------------------------------------
# Overview sheet: Status columns (zh-cn / de-de) now read "Handed back"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $handedBackStatus
$overview.Range("F2").Value = $handedBackStatus
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $handedBackStatus
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# Latest Target File -> hyperlink back to the source markdown file
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $mdUrl, "", "", $mdDisplay)
$zhcn.Range("I2").Font.Underline = $true
$zhcn.Range("I2").Font.Color = $hyperlinkColor

# Latest Handback File -> the generated zh-cn xliff
$zhcn.Range("J2").Value = "7eacbddf-b45d-4045-a9d6-3e4957d931f6.45c32c21ca1dbe28c7608e8c9b466d773d380610.zh-cn.xlf"

# Latest Handback DateTime -> real timestamp (was the zero-date placeholder)
$zhcn.Range("K2").Value = "2016-09-02 23:05:47"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $handedBackStatus
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

# Latest Target File -> hyperlink back to the source markdown file
$dede.Hyperlinks.Add($dede.Range("I2"), $mdUrl, "", "", $mdDisplay)
$dede.Range("I2").Font.Underline = $true
$dede.Range("I2").Font.Color = $hyperlinkColor

# Latest Handback File -> the generated de-de xliff
$dede.Range("J2").Value = "7eacbddf-b45d-4045-a9d6-3e4957d931f6.45c32c21ca1dbe28c7608e8c9b466d773d380610.de-de.xlf"

# Latest Handback DateTime -> real timestamp (de-de finished after zh-cn)
$dede.Range("K2").Value = "2016-09-02 23:05:55"
